$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Update "Recommandations" sheet (rows 2-47) ---
$ws1.Range("D2").Value = 3321.8
$ws1.Range("E2").Value = 101.59
$ws1.Range("D3").Value = 2780
$ws1.Range("E3").Value = 690
$ws1.Range("D4").Value = 2720
$ws1.Range("D5").Value = 2654.98
$ws1.Range("E5").Value = 663.97
$ws1.Range("D6").Value = 2360
$ws1.Range("A7").Value = "SETAO CI"
$ws1.Range("D7").Value = 2280
$ws1.Range("E7").Value = 590
$ws1.Range("A8").Value = "UNIWAX CI"
$ws1.Range("D8").Value = 2260
$ws1.Range("E8").Value = 565
$ws1.Range("D9").Value = 2090
$ws1.Range("E9").Value = 530
$ws1.Range("A10").Value = "BRVM - DISTRIBUTION"
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 1472.11
$ws1.Range("E10").Value = 362.37
$ws1.Range("A11").Value = "BRVM - TRANSPORT"
$ws1.Range("D11").Value = 1400.11
$ws1.Range("E11").Value = 353.7
$ws1.Range("A12").Value = "BRVM - AGRICULTURE"
$ws1.Range("D12").Value = 1221.66
$ws1.Range("E12").Value = 301.87
$ws1.Range("A13").Value = "SUCRIVOIRE"
$ws1.Range("C13").Value = 1
$ws1.Range("D13").Value = 995
$ws1.Range("E13").Value = 995
$ws1.Range("D14").Value = 831.74
$ws1.Range("E14").Value = 209.89
$ws1.Range("A15").Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Range("D15").Value = 709.47
$ws1.Range("E15").Value = 178.24
$ws1.Range("A16").Value = "BRVM-PRINCIPAL"
$ws1.Range("D16").Value = 707.99
$ws1.Range("E16").Value = 176.19
$ws1.Range("D17").Value = 525.76
$ws1.Range("E17").Value = 131.93
$ws1.Range("D18").Value = 522.14
$ws1.Range("E18").Value = 130.47
$ws1.Range("D19").Value = 489.77
$ws1.Range("E19").Value = 122.27
$ws1.Range("D20").Value = 481.34
$ws1.Range("E20").Value = 120.17
$ws1.Range("D21").Value = 430.92
$ws1.Range("E21").Value = 106.22
$ws1.Range("D22").Value = 429.39
$ws1.Range("E22").Value = 106.22
$ws1.Range("D23").Value = 376.45
$ws1.Range("E23").Value = 93.93
$ws1.Range("B24").Value = 3
$ws1.Range("D24").Value = 22.45
$ws1.Range("E24").Value = 7.47
$ws1.Range("F24").Value = "🟢 Achat"
$ws1.Range("G24").Value = "✅ Renforcer"
$ws1.Range("A25").Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Range("D25").Value = 8.29
$ws1.Range("E25").Value = 2.3
$ws1.Range("A27").Value = "SAFCA CI (SAFC)"
$ws1.Range("D27").Value = 6.92
$ws1.Range("E27").Value = 6.92
$ws1.Range("A28").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("D28").Value = 5.72
$ws1.Range("E28").Value = 5.72
$ws1.Range("A29").Value = "BICI CI (BICC)"
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 5.7
$ws1.Range("E29").Value = 5.7
$ws1.Range("G29").Value = "➖ Neutre"
$ws1.Range("A30").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Range("B30").Value = 2
$ws1.Range("D30").Value = 3.97
$ws1.Range("E30").Value = 1.54
$ws1.Range("A32").Value = "UNIWAX CI (UNXC)"
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = 2.12
$ws1.Range("E32").Value = 6.67
$ws1.Range("G32").Value = "👀 À surveiller"
$ws1.Range("A33").Value = "ONATEL BF (ONTBF)"
$ws1.Range("D33").Value = 2.01
$ws1.Range("E33").Value = 5.18
$ws1.Range("A34").Value = "SAPH CI (SPHC)"
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = 1.13
$ws1.Range("E34").Value = -2.13
$ws1.Range("G34").Value = "👀 À surveiller"
$ws1.Range("A35").Value = "AIR LIQUIDE CI (SIVC)"
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = 1.06
$ws1.Range("E35").Value = 3.92
$ws1.Range("A36").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("D36").Value = 0.37
$ws1.Range("E36").Value = -5.88
$ws1.Range("A37").Value = "TOTAL"
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 4
$ws1.Range("D37").Value = 0
$ws1.Range("E37").Value = 0
$ws1.Range("G37").Value = "➖ Neutre"
$ws1.Range("A38").Value = "SMB CI (SMBC)"
$ws1.Range("D38").Value = -0.09
$ws1.Range("E38").Value = -3.77
$ws1.Range("A39").Value = "BERNABE CI (BNBC)"
$ws1.Range("B39").Value = 2
$ws1.Range("C39").Value = 2
$ws1.Range("D39").Value = -1.01
$ws1.Range("E39").Value = 7.44
$ws1.Range("G39").Value = "👀 À surveiller"
$ws1.Range("A40").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("D40").Value = -2.13
$ws1.Range("E40").Value = -2.13
$ws1.Range("A41").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("D41").Value = -2.44
$ws1.Range("E41").Value = -2.44
$ws1.Range("A42").Value = "SOGB CI (SOGC)"
$ws1.Range("D42").Value = -3.04
$ws1.Range("E42").Value = -3.04
$ws1.Range("A43").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("B43").Value = 0
$ws1.Range("D43").Value = -3.78
$ws1.Range("E43").Value = -1.86
$ws1.Range("G43").Value = "➖ Neutre"
$ws1.Range("A44").Value = "SICOR CI (SICC)"
$ws1.Range("D44").Value = -5.71
$ws1.Range("E44").Value = -5.71
$ws1.Range("A45").Value = "SODE CI (SDCC)"
$ws1.Range("C45").Value = 2
$ws1.Range("D45").Value = -8.92
$ws1.Range("E45").Value = -6.61
$ws1.Range("A46").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("C46").Value = 2
$ws1.Range("D46").Value = -8.94
$ws1.Range("E46").Value = -5.9
$ws1.Range("A47").Value = "CIE CI (CIEC)"
$ws1.Range("C47").Value = 2
$ws1.Range("D47").Value = -10.25
$ws1.Range("E47").Value = -6.25

# --- Remove obsolete rows 48-50 (sheet now ends at row 47) ---
$ws1.Rows("48:50").Delete()

# --- Update "Top_YTD" sheet (rows 2-11) ---
$ws2.Range("B2").Value = 8556440.06
$ws2.Range("B3").Value = 399339.8
$ws2.Range("B4").Value = 370004.74
$ws2.Range("B5").Value = 340145.08
$ws2.Range("B6").Value = 226571.21
$ws2.Range("A7").Value = "SETAO CI"
$ws2.Range("B7").Value = 200828
$ws2.Range("A8").Value = "UNIWAX CI"
$ws2.Range("B8").Value = 195385.06
$ws2.Range("B9").Value = 150036.88
$ws2.Range("B10").Value = 47875.88
$ws2.Range("B11").Value = 40910.5

Write-Output "Recommandations and Top_YTD sheets updated successfully."
